# Rename "Sheet1" to "ADC" and move "TIM Generator" to sit right after it,
# so the tab order becomes: ADC, TIM Generator, Java processing, Triangle.
# The screen is not cleared/blinked while doing this (no Application.ScreenUpdating toggling).

$wb = $excel.ActiveWorkbook

$adc = $wb.Worksheets.Item("Sheet1")
$adc.Name = "ADC"

$tim = $wb.Worksheets.Item("TIM Generator")
$tim.Move($adc.Next)

# Worksheet handles in this host track position, not identity, so after the
# Move() re-resolve the sheets we still need by their (now current) names.
$adc = $wb.Worksheets.Item("ADC")
$tim = $wb.Worksheets.Item("TIM Generator")

# The chart on ADC still has its cached series formulas pointing at the old
# sheet name; update them so they reference the renamed sheet.
$chart = $adc.ChartObjects(1).Chart
$chart.SeriesCollection(1).Formula = "=SERIES(,ADC!`$E`$6:`$L`$6,ADC!`$E`$7:`$L`$7,1)"
$chart.SeriesCollection(2).Formula = "=SERIES(,ADC!`$E`$6:`$L`$6,ADC!`$E`$8:`$L`$8,2)"
$chart.SeriesCollection(3).Formula = "=SERIES(,ADC!`$E`$6:`$L`$6,ADC!`$E`$9:`$L`$9,3)"

# Restore per-sheet selections. Select TIM Generator's cell first, then
# finish on ADC so ADC remains the active/selected tab, matching the source.
$tim.Range("E15").Select() | Out-Null
$adc.Select() | Out-Null
$adc.Range("I4").Select() | Out-Null
